# Update ticket/view counts (column F) on the "展览" and "全部类型" sheets.
# Each entry maps a row number to its new value (column F) for a given sheet.

$wb = $excel.ActiveWorkbook

$sheet1Updates = @{
    2  = 189
    5  = 968
    6  = 5282
    7  = 444
    8  = 629
    12 = 29
    13 = 568
    14 = 15
    17 = 1762
    19 = 830
    21 = 188
    23 = 512
    24 = 131
    28 = 2594
    29 = 172
    31 = 51
    32 = 91
    33 = 23
    34 = 254
    40 = 641
    44 = 59
}

$sheet4Updates = @{
    3  = 189
    4  = 109
    5  = 968
    7  = 5282
    8  = 444
    9  = 629
    17 = 29
    18 = 568
    19 = 15
    23 = 1762
    25 = 830
    26 = 188
    29 = 512
    30 = 131
    33 = 2594
    34 = 172
    36 = 91
    37 = 23
    38 = 254
    46 = 59
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
